$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header info (row 2-3) ---
$ws.Range("C2").Value = "Hartmut"

# B3 holds a 16-digit account number that must stay TEXT (it was stored as
# inline string before the edit too). A direct Value assignment gets
# auto-detected as a number, so stage the apostrophe-forced text on a
# scratch cell (B9, which is going to be cleared out later anyway), move
# just the VALUE across (this preserves B3's own original style/format),
# then restore the scratch cell's format+value from its neighbour B10 so it
# ends up back in its normal blank state.
$ws.Range("B9").Value = "'2570314725427075"
$ws.Range("B9").Copy()
$ws.Range("B3").PasteSpecial(-4163)
$ws.Range("B10").Copy()
$ws.Range("B9").PasteSpecial(-4122)
$ws.Range("B9").Value = ""

$ws.Range("C3").Value = "Mohaupt"

# --- Opening balance line (row 5) ---
$ws.Range("D5").Value = "KONTOSTAND AM 26.12.2023"

# --- Row 6 (transaction) ---
$ws.Range("B6").Value = "27.12."
$ws.Range("C6").Value = "28.12."
$ws.Range("D6").Value = "KARTENZ./27.12 LIDL RO"
$ws.Range("E6").Value = "16,39-"

# --- Row 7 (transaction) ---
$ws.Range("B7").Value = "31.12."
$ws.Range("C7").Value = "01.01."
$ws.Range("D7").Value = "MCDONALDS Brilon"
$ws.Range("E7").Value = "18,80-"

# --- Row 8 (transaction) ---
$ws.Range("B8").Value = "02.01."
$ws.Range("C8").Value = "03.01."
$ws.Range("D8").Value = "BEITRAG Allianz SE K-99371477"
$ws.Range("E8").Value = "56,99-"

# --- Row 9: previously held a transaction (AMAZON.DE), now cleared out ---
# E9's style changes from s=17 (right aligned) to s=13 (center/center/wrap),
# matching the style already used by the neighbouring empty rows' B/C/D
# cells. Borrow that exact format via a format-only paste (keeps the style
# table free of duplicate/orphan entries), then clear the value and nudge
# the horizontal alignment to land precisely on the s=13 definition.
$ws.Range("B9").Copy()
$ws.Range("E9").PasteSpecial(-4122)
$ws.Range("E9").Value = ""
$ws.Range("E9").HorizontalAlignment = -4108

# B9/C9/D9 keep their existing style (s=8) but lose their content.
$ws.Range("B9").Value = ""
$ws.Range("C9").Value = ""
$ws.Range("D9").Value = ""

# --- Closing balance line (row 12) ---
$ws.Range("D12").Value = "KONTOSTAND AM 06.01.2024"
$ws.Range("E12").Value = "92,18-"

# --- Next billing date (row 13) ---
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 13.01.2024"
